# Update column F ("dSF") values for the rows that were re-pulled/recomputed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = 2
    9  = 1
    14 = 2
    16 = 1
    19 = 3
    22 = -2
    23 = 1
    24 = 0
    28 = 1
    29 = 1
    33 = -2
    36 = 2
    41 = -1
    43 = 7
    50 = 2
    63 = 0
    72 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
